# Add new data integration: new period_lbl rows for Apr/May/Jun 2024
# across the 2MATs / 3MMT / MAT / Month period types, rename the
# "Monthly" time_period_type to "Month", and highlight the newly added
# rows with a green fill.

$wb = $excel.ActiveWorkbook

$wsPeriod = $wb.Worksheets.Item("period_lbl")
$wsType   = $wb.Worksheets.Item("time_period_type")

# --- period_lbl: new rows 102-111 -----------------------------------
$newRows = @(
    @(101, "2MATs: Jun 2024", "2MATs: 2024 (06) Jun"),
    @(102, "3MMT: Apr 2024",  "3MMT: 2024 (04) Apr"),
    @(103, "3MMT: May 2024",  "3MMT: 2024 (05) May"),
    @(104, "3MMT: Jun 2024",  "3MMT: 2024 (06) Jun"),
    @(105, "MAT: Apr 2024",   "MAT: 2024 (04) Apr"),
    @(106, "MAT: May 2024",   "MAT: 2024 (05) May"),
    @(107, "MAT: Jun 2024",   "MAT: 2024 (06) Jun"),
    @(108, "Month: Apr 2024", "Month: 2024 (04) Apr"),
    @(109, "Month: May 2024", "Month: 2024 (05) May"),
    @(110, "Month: Jun 2024", "Month: 2024 (06) Jun")
)

$startRow = 102
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $code = $newRows[$i][0]
    $lbl  = $newRows[$i][1]
    $num  = $newRows[$i][2]

    $wsPeriod.Cells.Item($r, 1).Value = $code
    $wsPeriod.Cells.Item($r, 2).Value = $lbl
    $wsPeriod.Cells.Item($r, 3).Value = $num
}

# Highlight the new rows with a solid green fill (RGB 146,208,80 == FF92D050)
$wsPeriod.Range("A102:C111").Interior.Color = 5296274

# --- time_period_type: rename "Monthly" -> "Month" -------------------
$wsType.Range("B5").Value = "Month"

# --- header cells lose their old Consolas/centre-aligned styling -----
$wsPeriod.Range("B1").ClearFormats()
$wsType.Range("B1").ClearFormats()

# --- restore the cursor/selection to where the editor left it --------
$wsPeriod.Activate() | Out-Null
$wsPeriod.Range("C102").Select() | Out-Null

$wsType.Activate() | Out-Null
$wsType.Range("B6").Select() | Out-Null
